#
# Auto refresh - 16-02-2026 13:02:31.38
#
# Rebuilds the ML prediction columns on Excel_vs_ML:
#   - renames T/U/V headers
#   - adds two new columns (W: Risk_Level, X: Early_Warning)
#   - refreshes Feature_Importance ranking/weights
#   - refreshes Exec_Summary headline metrics
#

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Excel_vs_ML
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Excel_vs_ML")

# Re-label the existing ML columns and add the two new ones.
$ws.Range("T1").Value = "Predicted_Final_Deviation_%"
$ws.Range("U1").Value = "Risk_Score"
$ws.Range("V1").Value = "Predicted_Impact_Amount"
$ws.Range("W1").Value = "Risk_Level"
$ws.Range("X1").Value = "Early_Warning"

# Per-row model output: Predicted_Final_Deviation_%, Risk_Score,
# Predicted_Impact_Amount, Risk_Level, Early_Warning.
$rowData = @{
    2 = @($null, $null, $null, "LOW – Stable", "NO")
    3 = @($null, $null, $null, "LOW – Stable", "NO")
    4 = @($null, $null, $null, "LOW – Stable", "NO")
    5 = @(-18.71, 93.55, -95252.11792700001, "CRITICAL – Immediate Action", "NO")
    6 = @($null, $null, $null, "LOW – Stable", "NO")
    7 = @($null, $null, $null, "LOW – Stable", "NO")
    8 = @($null, $null, $null, "LOW – Stable", "NO")
    9 = @($null, $null, $null, "LOW – Stable", "NO")
    10 = @($null, $null, $null, "LOW – Stable", "NO")
    11 = @($null, $null, $null, "LOW – Stable", "NO")
    12 = @($null, $null, $null, "LOW – Stable", "NO")
    13 = @(-1.51, 7.55, -2440.23701, "LOW – Stable", "NO")
    14 = @($null, $null, $null, "LOW – Stable", "NO")
    15 = @(9.02, 45.09999999999999, 38215.927882, "MODERATE – Monitor Closely", "NO")
    16 = @(-22.13, 100, -23446.971791, "CRITICAL – Immediate Action", "NO")
    17 = @($null, $null, $null, "LOW – Stable", "NO")
    18 = @(-18.92, 94.60000000000001, -54085.323116, "CRITICAL – Immediate Action", "NO")
    19 = @($null, $null, $null, "LOW – Stable", "NO")
    20 = @($null, $null, $null, "LOW – Stable", "NO")
    21 = @(-0.19, 0.95, -373.876984, "LOW – Stable", "NO")
    22 = @($null, $null, $null, "LOW – Stable", "NO")
    23 = @(9.6, 48, 34143.38496, "MODERATE – Monitor Closely", "NO")
    24 = @($null, $null, $null, "LOW – Stable", "NO")
    25 = @($null, $null, $null, "LOW – Stable", "NO")
    26 = @($null, $null, $null, "LOW – Stable", "NO")
    27 = @(-18.44, 92.2, -109155.670356, "CRITICAL – Immediate Action", "NO")
    28 = @(-18.95, 94.75, -73977.35299499999, "CRITICAL – Immediate Action", "NO")
    29 = @(-1.6, 8, -2052.30064, "LOW – Stable", "NO")
    30 = @(-2.91, 14.55, -15540.11295, "LOW – Stable", "NO")
    31 = @($null, $null, $null, "LOW – Stable", "NO")
    32 = @(9.08, 45.4, 54151.20684400001, "MODERATE – Monitor Closely", "NO")
    33 = @($null, $null, $null, "LOW – Stable", "NO")
    34 = @(-16.29, 81.45, -93548.164347, "CRITICAL – Immediate Action", "NO")
    35 = @($null, $null, $null, "LOW – Stable", "NO")
    36 = @(-20.61, 100, -73085.329161, "CRITICAL – Immediate Action", "NO")
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]

    if ($vals[0] -eq $null) {
        $ws.Cells.Item($r, 20).Value = ""
    } else {
        $ws.Cells.Item($r, 20).Value = $vals[0]
    }

    if ($vals[1] -eq $null) {
        $ws.Cells.Item($r, 21).Value = ""
    } else {
        $ws.Cells.Item($r, 21).Value = $vals[1]
    }

    if ($vals[2] -eq $null) {
        $ws.Cells.Item($r, 22).Value = ""
    } else {
        $ws.Cells.Item($r, 22).Value = $vals[2]
    }

    $ws.Cells.Item($r, 23).Value = $vals[3]
    $ws.Cells.Item($r, 24).Value = $vals[4]
}

# ---------------------------------------------------------------------
# Sheet 2: Feature_Importance  (re-ranked after model retrain)
# ---------------------------------------------------------------------
$fi = $wb.Worksheets.Item("Feature_Importance")

$fi.Range("A2").Value = "Pace_Ratio"
$fi.Range("B2").Value = 0.8767043416575715

$fi.Range("A3").Value = "DSP_enc"
$fi.Range("B3").Value = 0.03902579892808492

$fi.Range("A4").Value = "Spend_Velocity"
$fi.Range("B4").Value = 0.03885957600547516

$fi.Range("A5").Value = "Spend_to_Date"
$fi.Range("B5").Value = 0.02013960227352946

$fi.Range("A6").Value = "Total_Budget"
$fi.Range("B6").Value = 0.0131773290904567

$fi.Range("A7").Value = "Flight_Days"
$fi.Range("B7").Value = 0.006460673297380543

$fi.Range("A8").Value = "Days_Elapsed"
$fi.Range("B8").Value = 0.005632678747501668

# ---------------------------------------------------------------------
# Sheet 3: Exec_Summary
# ---------------------------------------------------------------------
$es = $wb.Worksheets.Item("Exec_Summary")

$es.Range("A2").Value = "ML Validation MAE"
$es.Range("B2").Value = 0.045

$es.Range("A3").Value = "Total Predicted Impact"
$es.Range("B3").Value = -416446.94

$es.Range("A4").Value = "LAST_REFRESH_UTC"
$es.Range("B4").Value = "2026-02-16 07:32 UTC"
